$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 216
$ws.Range("I5").Value = 216
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 216
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -101
$ws.Range("N5").Value = ""
$ws.Range("H9").Value = 2679.4285
$ws.Range("I9").Value = 2743
$ws.Range("J9").Value = 2594.6667
$ws.Range("K9").Value = 2743
$ws.Range("L9").Value = 2594.6667
$ws.Range("M9").Value = -2574
$ws.Range("N9").Value = -2932.6667
$ws.Range("H12").Value = 4995
$ws.Range("I12").Value = 4990
$ws.Range("K12").Value = 4990
$ws.Range("M12").Value = -4820
$ws.Range("H19").Value = 54240.375
$ws.Range("I19").Value = 3785.3333
$ws.Range("J19").Value = 84513.39999999999
$ws.Range("K19").Value = 3785.3333
$ws.Range("L19").Value = 84513.39999999999
$ws.Range("M19").Value = -3610.3333
$ws.Range("N19").Value = -84863.39999999999
$ws.Range("H32").Value = 8749.5
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 11332.667
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 11332.667
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -11984.667
$ws.Range("H39").Value = 221.5
$ws.Range("I39").Value = 221.5
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 664.5
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -368.5
$ws.Range("N39").Value = ""
$ws.Range("H41").Value = 717.2353000000001
$ws.Range("I41").Value = 669.3333
$ws.Range("J41").Value = 771.125
$ws.Range("K41").Value = 669.3333
$ws.Range("L41").Value = 771.125
$ws.Range("M41").Value = -229.3333
$ws.Range("N41").Value = -1651.125
$ws.Range("H43").Value = 4570.603
$ws.Range("I43").Value = 4227.2183
$ws.Range("J43").Value = 6023.385
$ws.Range("K43").Value = 4227.2183
$ws.Range("L43").Value = 6023.385
$ws.Range("M43").Value = -4158.2183
$ws.Range("N43").Value = -6161.385
$ws.Range("H51").Value = 4928
$ws.Range("I51").Value = 4675.857
$ws.Range("J51").Value = 5088.4546
$ws.Range("K51").Value = 4675.857
$ws.Range("L51").Value = 5088.4546
$ws.Range("M51").Value = -4191.857
$ws.Range("N51").Value = -6056.4546
$ws.Range("H52").Value = 456.25
$ws.Range("H55").Value = 189
$ws.Range("J55").Value = 401
$ws.Range("L55").Value = 401
$ws.Range("N55").Value = -829
$ws.Range("H62").Value = 10333.167
$ws.Range("I62").Value = 7000
$ws.Range("K62").Value = 7000
$ws.Range("M62").Value = -6376
$ws.Range("H64").Value = 71433630
$ws.Range("J64").Value = 6233.3335
$ws.Range("L64").Value = 6233.3335
$ws.Range("N64").Value = -6729.3335
$ws.Range("H65").Value = 10333.167
$ws.Range("I65").Value = 7000
$ws.Range("K65").Value = 35000
$ws.Range("M65").Value = -31880
$ws.Range("H67").Value = 71433630
$ws.Range("J67").Value = 6233.3335
$ws.Range("L67").Value = 6233.3335
$ws.Range("N67").Value = -7949.3335
$ws.Range("H70").Value = 3419.889
$ws.Range("I70").Value = 2473.75
$ws.Range("J70").Value = 4176.8
$ws.Range("K70").Value = 7421.25
$ws.Range("L70").Value = 12530.4
$ws.Range("M70").Value = -7151.25
$ws.Range("N70").Value = -13070.4
$ws.Range("H73").Value = 3419.889
$ws.Range("I73").Value = 2473.75
$ws.Range("J73").Value = 4176.8
$ws.Range("K73").Value = 7421.25
$ws.Range("L73").Value = 12530.4
$ws.Range("M73").Value = -6485.25
$ws.Range("N73").Value = -14402.4
$ws.Range("H103").Value = 689.5
$ws.Range("J103").Value = 627
$ws.Range("L103").Value = 1881
$ws.Range("N103").Value = -3053
$ws.Range("H107").Value = 2078.2856
$ws.Range("I107").Value = 2078.2856
$ws.Range("K107").Value = 2078.2856
$ws.Range("M107").Value = -158.2856000000002
$ws.Range("H108").Value = 100000
$ws.Range("J108").Value = 100000
$ws.Range("L108").Value = 100000
$ws.Range("N108").Value = -107680
$ws.Range("H112").Value = 2523.5667
$ws.Range("J112").Value = 2888.28
$ws.Range("L112").Value = 8664.84
$ws.Range("N112").Value = -10880.84
$ws.Range("H113").Value = 22240138
$ws.Range("I113").Value = 40018350
$ws.Range("J113").Value = 17374.5
$ws.Range("K113").Value = 40018350
$ws.Range("L113").Value = 17374.5
$ws.Range("M113").Value = -40015096
$ws.Range("N113").Value = -23882.5
$ws.Range("H116").Value = 3500
$ws.Range("I116").Value = 3500
$ws.Range("K116").Value = 3500
$ws.Range("M116").Value = -58
$ws.Range("H127").Value = 2549999.5
$ws.Range("I127").Value = 2549999.5
$ws.Range("K127").Value = 7649998.5
$ws.Range("M127").Value = -7645038.5
$ws.Range("H132").Value = 8836.799999999999
$ws.Range("I132").Value = 6538.6055
$ws.Range("K132").Value = 19615.8165
$ws.Range("M132").Value = -17085.8165
$ws.Range("H135").Value = 2967.8708
$ws.Range("I135").Value = 1400.1333
$ws.Range("K135").Value = 12601.1997
$ws.Range("M135").Value = -10066.1997
$ws.Range("H137").Value = 10802.131
$ws.Range("I137").Value = 1787.4828
$ws.Range("J137").Value = 26180.059
$ws.Range("K137").Value = 5362.4484
$ws.Range("L137").Value = 78540.177
$ws.Range("M137").Value = -2812.4484
$ws.Range("N137").Value = -83640.177
$ws.Range("H138").Value = 3539.0386
$ws.Range("I138").Value = 4706.1816
$ws.Range("J138").Value = 3347.418
$ws.Range("K138").Value = 14118.5448
$ws.Range("L138").Value = 10042.254
$ws.Range("M138").Value = -8978.5448
$ws.Range("N138").Value = -20322.254
$ws.Range("H140").Value = 235153.22
$ws.Range("J140").Value = 235153.22
$ws.Range("L140").Value = 235153.22
$ws.Range("N140").Value = -245513.22
$ws.Range("H141").Value = 1895.5
$ws.Range("I141").Value = 999.5
$ws.Range("J141").Value = 2343.5
$ws.Range("K141").Value = 2998.5
$ws.Range("L141").Value = 7030.5
$ws.Range("M141").Value = 2181.5
$ws.Range("N141").Value = -17390.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1080
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -1626
$ws.Range("H32").Value = 5325.5835
$ws.Range("I32").Value = 2991.3594
$ws.Range("K32").Value = 2991.3594
$ws.Range("M32").Value = -2704.3594
$ws.Range("H33").Value = 2023.8096
$ws.Range("I33").Value = 2023.8096
$ws.Range("K33").Value = 2023.8096
$ws.Range("M33").Value = -1694.8096
$ws.Range("H36").Value = 2516.3333
$ws.Range("I36").Value = 2516.3333
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2516.3333
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2170.3333
$ws.Range("N36").Value = ""
$ws.Range("H37").Value = 12017
$ws.Range("I37").Value = 12017
$ws.Range("K37").Value = 12017
$ws.Range("M37").Value = -11744
$ws.Range("H45").Value = 1568.6
$ws.Range("I45").Value = 1613.6666
$ws.Range("K45").Value = 1613.6666
$ws.Range("M45").Value = -1236.6666
$ws.Range("H61").Value = 9883.048000000001
$ws.Range("I61").Value = 2277.5925
$ws.Range("K61").Value = 2277.5925
$ws.Range("M61").Value = -2065.5925
$ws.Range("H74").Value = 12690.18
$ws.Range("I74").Value = 2347.6667
$ws.Range("K74").Value = 2347.6667
$ws.Range("M74").Value = -1473.6667
$ws.Range("H77").Value = 12690.18
$ws.Range("I77").Value = 2347.6667
$ws.Range("K77").Value = 11738.3335
$ws.Range("M77").Value = -7370.333500000001
$ws.Range("H88").Value = 2105.1428
$ws.Range("J88").Value = 2224.182
$ws.Range("L88").Value = 2224.182
$ws.Range("N88").Value = -3036.182
$ws.Range("H91").Value = 2105.1428
$ws.Range("J91").Value = 2224.182
$ws.Range("L91").Value = 2224.182
$ws.Range("N91").Value = -5032.182
$ws.Range("H110").Value = 4768.391
$ws.Range("I110").Value = 6255
$ws.Range("J110").Value = 1981
$ws.Range("K110").Value = 6255
$ws.Range("L110").Value = 1981
$ws.Range("M110").Value = -4210
$ws.Range("N110").Value = -6071
$ws.Range("H116").Value = 1080
$ws.Range("I116").Value = 1000
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 1000
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = 1294
$ws.Range("N116").Value = -5988
$ws.Range("H122").Value = 1404416.5
$ws.Range("H126").Value = 6442.778
$ws.Range("I126").Value = 6442.778
$ws.Range("K126").Value = 19328.334
$ws.Range("M126").Value = -16858.334
$ws.Range("H132").Value = 2231876.8
$ws.Range("I132").Value = 3146.3103
$ws.Range("K132").Value = 9438.930899999999
$ws.Range("M132").Value = -6908.930899999999
$ws.Range("H136").Value = 9883.048000000001
$ws.Range("I136").Value = 2277.5925
$ws.Range("K136").Value = 6832.7775
$ws.Range("M136").Value = -4282.7775

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1080
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1400
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = -886
$ws.Range("N3").Value = -1628
$ws.Range("H20").Value = 8355164
$ws.Range("I20").Value = 20840942
$ws.Range("J20").Value = 31312.334
$ws.Range("K20").Value = 20840942
$ws.Range("L20").Value = 31312.334
$ws.Range("M20").Value = -20840695
$ws.Range("N20").Value = -31806.334
$ws.Range("H22").Value = 10205741
$ws.Range("I22").Value = 10205741
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 10205741
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -10205568
$ws.Range("N22").Value = ""
$ws.Range("H74").Value = 59992.5
$ws.Range("J74").Value = 59992.5
$ws.Range("L74").Value = 59992.5
$ws.Range("N74").Value = -61864.5
$ws.Range("H75").Value = 6665.2856
$ws.Range("I75").Value = 6665.2856
$ws.Range("K75").Value = 6665.2856
$ws.Range("M75").Value = -5729.2856
$ws.Range("H76").Value = 24000
$ws.Range("J76").Value = 24000
$ws.Range("L76").Value = 24000
$ws.Range("N76").Value = -24630
$ws.Range("H77").Value = 59992.5
$ws.Range("J77").Value = 59992.5
$ws.Range("L77").Value = 179977.5
$ws.Range("N77").Value = -189337.5
$ws.Range("H78").Value = 6665.2856
$ws.Range("I78").Value = 6665.2856
$ws.Range("K78").Value = 19995.8568
$ws.Range("M78").Value = -15315.8568
$ws.Range("H79").Value = 24000
$ws.Range("J79").Value = 24000
$ws.Range("L79").Value = 24000
$ws.Range("N79").Value = -26184
$ws.Range("H94").Value = 5061.524
$ws.Range("I94").Value = 3362
$ws.Range("J94").Value = 10500
$ws.Range("K94").Value = 3362
$ws.Range("L94").Value = 10500
$ws.Range("M94").Value = -2911
$ws.Range("N94").Value = -11402
$ws.Range("H134").Value = 12029.233
$ws.Range("I134").Value = 5836.5713
$ws.Range("K134").Value = 17509.7139
$ws.Range("M134").Value = -14974.7139
$ws.Range("H138").Value = 160780
$ws.Range("J138").Value = 160780
$ws.Range("L138").Value = 160780
$ws.Range("N138").Value = -171060

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 62.5
$ws.Range("I17").Value = 62.5
$ws.Range("K17").Value = 62.5
$ws.Range("M17").Value = 111.5
$ws.Range("H22").Value = 1228.0294
$ws.Range("I22").Value = 806.4400000000001
$ws.Range("J22").Value = 2399.111
$ws.Range("K22").Value = 806.4400000000001
$ws.Range("L22").Value = 2399.111
$ws.Range("M22").Value = -456.4400000000001
$ws.Range("N22").Value = -3099.111
$ws.Range("H31").Value = 67128.59
$ws.Range("I31").Value = 101918.336
$ws.Range("J31").Value = 21467.062
$ws.Range("K31").Value = 101918.336
$ws.Range("L31").Value = 21467.062
$ws.Range("M31").Value = -101623.336
$ws.Range("N31").Value = -22057.062
$ws.Range("H34").Value = 67128.59
$ws.Range("I34").Value = 101918.336
$ws.Range("J34").Value = 21467.062
$ws.Range("K34").Value = 101918.336
$ws.Range("L34").Value = 21467.062
$ws.Range("M34").Value = -101716.336
$ws.Range("N34").Value = -21871.062
$ws.Range("H58").Value = 12956.275
$ws.Range("I58").Value = 4277.0454
$ws.Range("J58").Value = 23564.223
$ws.Range("K58").Value = 4277.0454
$ws.Range("L58").Value = 23564.223
$ws.Range("M58").Value = -4074.0454
$ws.Range("N58").Value = -23970.223
$ws.Range("H74").Value = 38500
$ws.Range("J74").Value = 38500
$ws.Range("L74").Value = 38500
$ws.Range("N74").Value = -40248
$ws.Range("H77").Value = 38500
$ws.Range("J77").Value = 38500
$ws.Range("L77").Value = 115500
$ws.Range("N77").Value = -124236
$ws.Range("H94").Value = 1351
$ws.Range("I94").Value = 510.66666
$ws.Range("K94").Value = 510.66666
$ws.Range("M94").Value = -59.66665999999998
$ws.Range("H99").Value = 6306358.5
$ws.Range("I99").Value = 8904951
$ws.Range("K99").Value = 8904951
$ws.Range("M99").Value = -8903453
$ws.Range("H103").Value = 10537.833
$ws.Range("I103").Value = 5931.75
$ws.Range("K103").Value = 5931.75
$ws.Range("M103").Value = -4759.75
$ws.Range("H107").Value = 1147.1351
$ws.Range("I107").Value = 937.92
$ws.Range("J107").Value = 1583
$ws.Range("K107").Value = 937.92
$ws.Range("L107").Value = 1583
$ws.Range("M107").Value = 982.08
$ws.Range("N107").Value = -5423
$ws.Range("H122").Value = 1994.2084
$ws.Range("I122").Value = 2817.889
$ws.Range("K122").Value = 8453.667000000001
$ws.Range("M122").Value = -6003.667000000001
$ws.Range("H126").Value = 6306358.5
$ws.Range("I126").Value = 8904951
$ws.Range("K126").Value = 26714853
$ws.Range("M126").Value = -26712383
$ws.Range("H132").Value = 5276
$ws.Range("I132").Value = 2147.348
$ws.Range("J132").Value = 13271.444
$ws.Range("K132").Value = 6442.044
$ws.Range("L132").Value = 39814.33199999999
$ws.Range("M132").Value = -3912.044
$ws.Range("N132").Value = -44874.33199999999
$ws.Range("H134").Value = 37043660
$ws.Range("I134").Value = 1429.4
$ws.Range("K134").Value = 4288.200000000001
$ws.Range("M134").Value = -1753.200000000001
$ws.Range("H136").Value = 12956.275
$ws.Range("I136").Value = 4277.0454
$ws.Range("J136").Value = 23564.223
$ws.Range("K136").Value = 12831.1362
$ws.Range("L136").Value = 70692.66900000001
$ws.Range("M136").Value = -10281.1362
$ws.Range("N136").Value = -75792.66900000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 817.2069
$ws.Range("I7").Value = 127.041664
$ws.Range("J7").Value = 4130
$ws.Range("K7").Value = 381.124992
$ws.Range("L7").Value = 12390
$ws.Range("M7").Value = -269.124992
$ws.Range("N7").Value = -12614
$ws.Range("H10").Value = 275.44446
$ws.Range("I10").Value = 278.75
$ws.Range("K10").Value = 836.25
$ws.Range("M10").Value = -697.25
$ws.Range("H16").Value = 503.125
$ws.Range("I16").Value = 579.1667
$ws.Range("J16").Value = 275
$ws.Range("K16").Value = 1737.5001
$ws.Range("L16").Value = 825
$ws.Range("M16").Value = -1564.5001
$ws.Range("N16").Value = -1171
$ws.Range("H46").Value = 498.33334
$ws.Range("J46").Value = 999.3333
$ws.Range("L46").Value = 2997.9999
$ws.Range("N46").Value = -3179.9999
$ws.Range("H63").Value = 23951.5
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 23951.5
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 71854.5
$ws.Range("M63").Value = ""
$ws.Range("N63").Value = -73352.5
$ws.Range("H66").Value = 23951.5
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 23951.5
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 215563.5
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = -223051.5
$ws.Range("H80").Value = 16745.75
$ws.Range("J80").Value = 16745.75
$ws.Range("L80").Value = 50237.25
$ws.Range("N80").Value = -52109.25
$ws.Range("H83").Value = 16745.75
$ws.Range("J83").Value = 16745.75
$ws.Range("L83").Value = 150711.75
$ws.Range("N83").Value = -160071.75
$ws.Range("H121").Value = 9995
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 9995
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 29985
$ws.Range("M121").Value = ""
$ws.Range("N121").Value = -32605
$ws.Range("H122").Value = 6152223
$ws.Range("I122").Value = 16988656
$ws.Range("K122").Value = 152897904
$ws.Range("M122").Value = -152895454
$ws.Range("H131").Value = 1473.87
$ws.Range("J131").Value = 1473.87
$ws.Range("L131").Value = 4421.61
$ws.Range("N131").Value = -14501.61
$ws.Range("H132").Value = 1999.375
$ws.Range("J132").Value = 1819.8
$ws.Range("L132").Value = 16378.2
$ws.Range("N132").Value = -21438.2
$ws.Range("H134").Value = 7104.213
$ws.Range("I134").Value = 4771.143
$ws.Range("K134").Value = 14313.429
$ws.Range("M134").Value = -9243.429
$ws.Range("H137").Value = 25121
$ws.Range("I137").Value = 209
$ws.Range("K137").Value = 627
$ws.Range("M137").Value = 4473
$ws.Range("H138").Value = 4375.4287
$ws.Range("I138").Value = 3688
$ws.Range("K138").Value = 11064
$ws.Range("M138").Value = -5924
$ws.Range("H140").Value = 2188
$ws.Range("I140").Value = 2025.6
$ws.Range("K140").Value = 6076.799999999999
$ws.Range("M140").Value = -896.7999999999993

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13303.4
$ws.Range("I43").Value = 14172.333
$ws.Range("J43").Value = 12000
$ws.Range("K43").Value = 14172.333
$ws.Range("L43").Value = 12000
$ws.Range("M43").Value = -14021.333
$ws.Range("N43").Value = -12302
$ws.Range("H70").Value = 14601.591
$ws.Range("I70").Value = 16217.25
$ws.Range("K70").Value = 16217.25
$ws.Range("M70").Value = -15947.25
$ws.Range("H73").Value = 14601.591
$ws.Range("I73").Value = 16217.25
$ws.Range("K73").Value = 16217.25
$ws.Range("M73").Value = -15281.25
$ws.Range("H102").Value = 20844334
$ws.Range("I102").Value = 25012800
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 25012800
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -25011178
$ws.Range("N102").Value = -5244
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080
$ws.Range("H113").Value = 1998.6
$ws.Range("I113").Value = 1998.25
$ws.Range("K113").Value = 1998.25
$ws.Range("M113").Value = 171.75
$ws.Range("I122").Value = 2122276
$ws.Range("J122").Value = 2949.75
$ws.Range("K122").Value = 6366828
$ws.Range("L122").Value = 8849.25
$ws.Range("M122").Value = -6364378
$ws.Range("N122").Value = -13749.25
$ws.Range("H132").Value = 13656.786
$ws.Range("I132").Value = 6316.4546
$ws.Range("J132").Value = 40571.332
$ws.Range("K132").Value = 18949.3638
$ws.Range("L132").Value = 121713.996
$ws.Range("M132").Value = -16419.3638
$ws.Range("N132").Value = -126773.996
$ws.Range("H133").Value = 65480.125
$ws.Range("J133").Value = 65480.125
$ws.Range("L133").Value = 65480.125
$ws.Range("N133").Value = -75600.125
$ws.Range("H139").Value = 61777.2
$ws.Range("J139").Value = 61777.2
$ws.Range("L139").Value = 61777.2
$ws.Range("N139").Value = -72057.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1708468.4
$ws.Range("I7").Value = 2652951
$ws.Range("K7").Value = 2652951
$ws.Range("M7").Value = -2652839
$ws.Range("H16").Value = 22226368
$ws.Range("I16").Value = 4412.5
$ws.Range("J16").Value = 66670276
$ws.Range("K16").Value = 4412.5
$ws.Range("L16").Value = 66670276
$ws.Range("M16").Value = -4242.5
$ws.Range("N16").Value = -66670616
$ws.Range("H22").Value = 2603.0217
$ws.Range("I22").Value = 2346.8286
$ws.Range("J22").Value = 3418.182
$ws.Range("K22").Value = 2346.8286
$ws.Range("L22").Value = 3418.182
$ws.Range("M22").Value = -2051.8286
$ws.Range("N22").Value = -4008.182
$ws.Range("H27").Value = 2603.0217
$ws.Range("I27").Value = 2346.8286
$ws.Range("J27").Value = 3418.182
$ws.Range("K27").Value = 2346.8286
$ws.Range("L27").Value = 3418.182
$ws.Range("M27").Value = -2239.8286
$ws.Range("N27").Value = -3632.182
$ws.Range("H40").Value = 2945081
$ws.Range("I40").Value = 1674.1666
$ws.Range("K40").Value = 1674.1666
$ws.Range("M40").Value = -1538.1666
$ws.Range("H46").Value = 2926.75
$ws.Range("J46").Value = 3999.111
$ws.Range("L46").Value = 3999.111
$ws.Range("N46").Value = -4375.111
$ws.Range("H55").Value = 1395.6428
$ws.Range("I55").Value = 992.5
$ws.Range("J55").Value = 2121.3
$ws.Range("K55").Value = 992.5
$ws.Range("L55").Value = 2121.3
$ws.Range("M55").Value = -819.5
$ws.Range("N55").Value = -2467.3
$ws.Range("H68").Value = 3284.4
$ws.Range("I68").Value = 2834.8572
$ws.Range("K68").Value = 2834.8572
$ws.Range("M68").Value = -2085.8572
$ws.Range("H71").Value = 3284.4
$ws.Range("I71").Value = 2834.8572
$ws.Range("K71").Value = 14174.286
$ws.Range("M71").Value = -10430.286
$ws.Range("H93").Value = 4890.522
$ws.Range("J93").Value = 7323.7144
$ws.Range("L93").Value = 7323.7144
$ws.Range("N93").Value = -9819.714400000001
$ws.Range("H101").Value = 17812.4
$ws.Range("J101").Value = 17812.4
$ws.Range("L101").Value = 17812.4
$ws.Range("N101").Value = -24302.4
$ws.Range("H122").Value = 7606.316
$ws.Range("I122").Value = 6543.3335
$ws.Range("J122").Value = 9428.571
$ws.Range("K122").Value = 19630.0005
$ws.Range("L122").Value = 28285.713
$ws.Range("M122").Value = -17180.0005
$ws.Range("N122").Value = -33185.713
$ws.Range("H126").Value = 1708468.4
$ws.Range("I126").Value = 2652951
$ws.Range("K126").Value = 7958853
$ws.Range("M126").Value = -7956383
$ws.Range("H132").Value = 5038250
$ws.Range("I132").Value = 39999
$ws.Range("J132").Value = 6704334
$ws.Range("K132").Value = 119997
$ws.Range("L132").Value = 20113002
$ws.Range("M132").Value = -117467
$ws.Range("N132").Value = -20118062
$ws.Range("H136").Value = 19447.65
$ws.Range("I136").Value = 22114.9
$ws.Range("J136").Value = 16780.4
$ws.Range("K136").Value = 66344.70000000001
$ws.Range("L136").Value = 50341.2
$ws.Range("M136").Value = -63794.70000000001
$ws.Range("N136").Value = -55441.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 13625
$ws.Range("H29").Value = 18006.666
$ws.Range("J29").Value = 10010
$ws.Range("L29").Value = 10010
$ws.Range("N29").Value = -10590
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("H62").Value = 155491.42
$ws.Range("I62").Value = 2250
$ws.Range("K62").Value = 2250
$ws.Range("M62").Value = -1626
$ws.Range("H65").Value = 155491.42
$ws.Range("I65").Value = 2250
$ws.Range("K65").Value = 11250
$ws.Range("M65").Value = -8130
$ws.Range("H74").Value = 44444
$ws.Range("J74").Value = 44444
$ws.Range("L74").Value = 44444
$ws.Range("N74").Value = -46316
$ws.Range("H77").Value = 44444
$ws.Range("J77").Value = 44444
$ws.Range("L77").Value = 133332
$ws.Range("N77").Value = -142692
$ws.Range("H82").Value = 47500
$ws.Range("I82").Value = 30000
$ws.Range("J82").Value = 53333.332
$ws.Range("K82").Value = 30000
$ws.Range("L82").Value = 53333.332
$ws.Range("M82").Value = -29617
$ws.Range("N82").Value = -54099.332
$ws.Range("H85").Value = 47500
$ws.Range("I85").Value = 30000
$ws.Range("J85").Value = 53333.332
$ws.Range("K85").Value = 30000
$ws.Range("L85").Value = 53333.332
$ws.Range("M85").Value = -28674
$ws.Range("N85").Value = -55985.332
$ws.Range("H107").Value = 1102.6666
$ws.Range("I107").Value = 1456.6
$ws.Range("J107").Value = 394.8
$ws.Range("K107").Value = 4369.799999999999
$ws.Range("L107").Value = 1184.4
$ws.Range("M107").Value = -2449.799999999999
$ws.Range("N107").Value = -5024.4
$ws.Range("H113").Value = 4648.615
$ws.Range("I113").Value = 6303.5557
$ws.Range("K113").Value = 18910.6671
$ws.Range("M113").Value = -16740.6671
$ws.Range("H122").Value = 4479.524
$ws.Range("I122").Value = 1697.6923
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 5093.0769
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -2643.0769
$ws.Range("N122").Value = -31900
$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""
$ws.Range("H132").Value = 4181.213
$ws.Range("I132").Value = 1985.0385
$ws.Range("J132").Value = 16870.223
$ws.Range("K132").Value = 5955.1155
$ws.Range("L132").Value = 50610.66900000001
$ws.Range("M132").Value = -3425.1155
$ws.Range("N132").Value = -55670.66900000001
$ws.Range("H136").Value = 11103.117
$ws.Range("I136").Value = 1124.3334
$ws.Range("J136").Value = 18981.105
$ws.Range("K136").Value = 3373.0002
$ws.Range("L136").Value = 56943.315
$ws.Range("M136").Value = -823.0001999999999
$ws.Range("N136").Value = -62043.315
$ws.Range("H138").Value = 100428.5
$ws.Range("J138").Value = 100428.5
$ws.Range("L138").Value = 100428.5
$ws.Range("N138").Value = -110708.5
$ws.Range("H141").Value = 94989
$ws.Range("J141").Value = 94989
$ws.Range("L141").Value = 94989
$ws.Range("N141").Value = -105349
